# Auto-generated Excel COM-interop script to apply the diff
$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update Last Updated timestamp ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("A2").Value = "30 Oct 2025, 09:27 AM"

# --- Top Gainers sheet ---
$wsGainers = $wb.Worksheets.Item("Top Gainers")
$wsGainers.Range("C2").Value = 7.6469
$wsGainers.Range("D2").Value = 15.2357
$wsGainers.Range("E2").Value = 28.3638
$wsGainers.Range("C3").Value = 7.5034
$wsGainers.Range("D3").Value = 8.025700000000001
$wsGainers.Range("E3").Value = 7.545
$wsGainers.Range("C4").Value = 5.6937
$wsGainers.Range("D4").Value = 19.0875
$wsGainers.Range("E4").Value = 19.1763
$wsGainers.Range("C6").Value = 4.8207
$wsGainers.Range("D6").Value = 6.669
$wsGainers.Range("E6").Value = 3.9823
$wsGainers.Range("B7").Value = "IFCI"
$wsGainers.Range("C7").Value = 3.8549
$wsGainers.Range("D7").Value = 6.4667
$wsGainers.Range("E7").Value = 9.6577
$wsGainers.Range("B8").Value = "KELLTONTEC"
$wsGainers.Range("C8").Value = 3.8462
$wsGainers.Range("D8").Value = 1.0695
$wsGainers.Range("E8").Value = -12.9032
$wsGainers.Range("B9").Value = "POLICYBZR"
$wsGainers.Range("C9").Value = 3.7484
$wsGainers.Range("D9").Value = 6.0687
$wsGainers.Range("E9").Value = 5.0529
$wsGainers.Range("B10").Value = "PDSL"
$wsGainers.Range("C10").Value = 3.6188
$wsGainers.Range("D10").Value = 6.9017
$wsGainers.Range("E10").Value = 12.9499
$wsGainers.Range("B11").Value = "VAIBHAVGBL"
$wsGainers.Range("C11").Value = 3.4971
$wsGainers.Range("D11").Value = 11.6952
$wsGainers.Range("E11").Value = 18.3586
$wsGainers.Range("B12").Value = "SUBROS"
$wsGainers.Range("C12").Value = 3.3955
$wsGainers.Range("D12").Value = 4.9264
$wsGainers.Range("E12").Value = 5.7772
$wsGainers.Range("B13").Value = "MEGASOFT"
$wsGainers.Range("C13").Value = 3.1551
$wsGainers.Range("D13").Value = 19.4111
$wsGainers.Range("E13").Value = 37.74
$wsGainers.Range("B14").Value = "SHREEJISPG"
$wsGainers.Range("C14").Value = 3.0361
$wsGainers.Range("D14").Value = 7.1865
$wsGainers.Range("E14").Value = 8.249599999999999
$wsGainers.Range("B15").Value = "IIFL"
$wsGainers.Range("C15").Value = 2.9828
$wsGainers.Range("D15").Value = 9.963200000000001
$wsGainers.Range("E15").Value = 19.1857
$wsGainers.Range("B16").Value = "BLS"
$wsGainers.Range("C16").Value = 2.8493
$wsGainers.Range("D16").Value = -0.1885
$wsGainers.Range("E16").Value = -1.4428
$wsGainers.Range("B17").Value = "V2RETAIL"
$wsGainers.Range("C17").Value = 2.7042
$wsGainers.Range("D17").Value = 1.2191
$wsGainers.Range("E17").Value = 9.301399999999999
$wsGainers.Range("B18").Value = "REFEX"
$wsGainers.Range("C18").Value = 2.6798
$wsGainers.Range("D18").Value = 0.11
$wsGainers.Range("E18").Value = 2.0895
$wsGainers.Range("B19").Value = "TEXRAIL"
$wsGainers.Range("C19").Value = 2.6593
$wsGainers.Range("D19").Value = 4.1738
$wsGainers.Range("E19").Value = 2.1964
$wsGainers.Range("B20").Value = "SALASAR"
$wsGainers.Range("C20").Value = 2.6396
$wsGainers.Range("D20").Value = 7.5532
$wsGainers.Range("E20").Value = 13.9797
$wsGainers.Range("B21").Value = "BAYERCROP"
$wsGainers.Range("C21").Value = 2.6223
$wsGainers.Range("D21").Value = -0.755
$wsGainers.Range("E21").Value = 3.561
$wsGainers.Range("B22").Value = "SDBL"
$wsGainers.Range("C22").Value = 2.5915
$wsGainers.Range("D22").Value = 1.1363
$wsGainers.Range("E22").Value = 6.7278
$wsGainers.Range("B23").Value = "BUTTERFLY"
$wsGainers.Range("C23").Value = 2.5551
$wsGainers.Range("D23").Value = 13.1136
$wsGainers.Range("E23").Value = 15.7733
$wsGainers.Range("B24").Value = "DEEDEV"
$wsGainers.Range("C24").Value = 2.5309
$wsGainers.Range("D24").Value = -4.2903
$wsGainers.Range("E24").Value = -5.0797
$wsGainers.Range("B25").Value = "SHANTIGOLD"
$wsGainers.Range("C25").Value = 2.4992
$wsGainers.Range("D25").Value = 9.7135
$wsGainers.Range("E25").Value = 2.3938
$wsGainers.Range("B26").Value = "KRISHANA"
$wsGainers.Range("C26").Value = 2.4221
$wsGainers.Range("D26").Value = -0.6163999999999999
$wsGainers.Range("E26").Value = 14.4664
$wsGainers.Range("B27").Value = "EIEL"
$wsGainers.Range("C27").Value = 2.4154
$wsGainers.Range("D27").Value = 3.5023
$wsGainers.Range("E27").Value = 4.2942
$wsGainers.Range("B28").Value = "SEQUENT"
$wsGainers.Range("C28").Value = 2.3981
$wsGainers.Range("D28").Value = 7.7749
$wsGainers.Range("E28").Value = 16.8984
$wsGainers.Range("B29").Value = "RPOWER"
$wsGainers.Range("C29").Value = 2.2462
$wsGainers.Range("D29").Value = 4.7345
$wsGainers.Range("E29").Value = 6.7178
$wsGainers.Range("B30").Value = "GMMPFAUDLR"
$wsGainers.Range("C30").Value = 2.2311
$wsGainers.Range("D30").Value = 6.6616
$wsGainers.Range("E30").Value = 18.911
$wsGainers.Range("B31").Value = "RAMASTEEL"
$wsGainers.Range("C31").Value = 2.2088
$wsGainers.Range("D31").Value = 2.1063
$wsGainers.Range("E31").Value = 3.666
$wsGainers.Range("B32").Value = "PRICOLLTD"
$wsGainers.Range("C32").Value = 2.167
$wsGainers.Range("D32").Value = 5.2808
$wsGainers.Range("E32").Value = 2.6636
$wsGainers.Range("B33").Value = "UNIPARTS"
$wsGainers.Range("C33").Value = 2.1648
$wsGainers.Range("D33").Value = 4.3245
$wsGainers.Range("E33").Value = 19.6848
$wsGainers.Range("B34").Value = "INTELLECT"
$wsGainers.Range("C34").Value = 2.1528
$wsGainers.Range("D34").Value = 5.3209
$wsGainers.Range("E34").Value = 7.8443
$wsGainers.Range("B35").Value = "ROSSTECH"
$wsGainers.Range("C35").Value = 2.1394
$wsGainers.Range("D35").Value = 4.0829
$wsGainers.Range("E35").Value = -4.8119
$wsGainers.Range("B36").Value = "PRABHA"
$wsGainers.Range("C36").Value = 2.1039
$wsGainers.Range("D36").Value = 4.3582
$wsGainers.Range("E36").Value = 0.6327
$wsGainers.Range("B37").Value = "SKFINDIA"
$wsGainers.Range("C37").Value = 2.098
$wsGainers.Range("D37").Value = -0.4274
$wsGainers.Range("E37").Value = -7.1211
$wsGainers.Range("B38").Value = "AXISCADES"
$wsGainers.Range("C38").Value = 2.0605
$wsGainers.Range("D38").Value = 9.688599999999999
$wsGainers.Range("E38").Value = -0.5646
$wsGainers.Range("B39").Value = "ICIL"
$wsGainers.Range("C39").Value = 2.0154
$wsGainers.Range("D39").Value = 1.4111
$wsGainers.Range("E39").Value = 7.6886
$wsGainers.Range("B40").Value = "ASALCBR"
$wsGainers.Range("C40").Value = 2.0084
$wsGainers.Range("D40").Value = 2.199
$wsGainers.Range("E40").Value = 15.1836
$wsGainers.Range("B41").Value = "SIGNATURE"
$wsGainers.Range("C41").Value = 1.9461
$wsGainers.Range("D41").Value = 0.8132
$wsGainers.Range("E41").Value = 5.8634
$wsGainers.Range("B42").Value = "ENDURANCE"
$wsGainers.Range("C42").Value = 1.8942
$wsGainers.Range("D42").Value = -0.4438
$wsGainers.Range("E42").Value = 5.4127
$wsGainers.Range("B43").Value = "BLSE"
$wsGainers.Range("C43").Value = 1.8577
$wsGainers.Range("D43").Value = 0.8588
$wsGainers.Range("E43").Value = -5.1255
$wsGainers.Range("B44").Value = "DATAMATICS"
$wsGainers.Range("C44").Value = 1.8562
$wsGainers.Range("D44").Value = 9.3072
$wsGainers.Range("E44").Value = 17.878
$wsGainers.Range("B45").Value = "CEMPRO"
$wsGainers.Range("C45").Value = 1.82
$wsGainers.Range("D45").Value = 7.7409
$wsGainers.Range("E45").Value = 6.3314
$wsGainers.Range("B46").Value = "NETWEB"
$wsGainers.Range("C46").Value = 1.8016
$wsGainers.Range("D46").Value = 7.4736
$wsGainers.Range("E46").Value = 9.494300000000001
$wsGainers.Range("B47").Value = "DOLPHIN"
$wsGainers.Range("C47").Value = 1.7956
$wsGainers.Range("D47").Value = 0.7887
$wsGainers.Range("E47").Value = -3.2713
$wsGainers.Range("B48").Value = "ECLERX"
$wsGainers.Range("C48").Value = 1.7846
$wsGainers.Range("D48").Value = 6.2025
$wsGainers.Range("E48").Value = 17.5984
$wsGainers.Range("B49").Value = "HYUNDAI"
$wsGainers.Range("C49").Value = 1.7517
$wsGainers.Range("D49").Value = 5.8367
$wsGainers.Range("E49").Value = -7.1738
$wsGainers.Range("B50").Value = "BAJAJINDEF"
$wsGainers.Range("C50").Value = 1.7114
$wsGainers.Range("D50").Value = 2.8011
$wsGainers.Range("E50").Value = 9.772600000000001
$wsGainers.Range("B51").Value = "SPMLINFRA"
$wsGainers.Range("C51").Value = 1.6951
$wsGainers.Range("D51").Value = 3.4626
$wsGainers.Range("E51").Value = -3.026
$wsGainers.Range("B52").Value = "RAYMONDLSL"
$wsGainers.Range("C52").Value = 1.6896
$wsGainers.Range("D52").Value = 5.715
$wsGainers.Range("E52").Value = 6.2613
$wsGainers.Range("B53").Value = "GMRAIRPORT"
$wsGainers.Range("C53").Value = 1.678
$wsGainers.Range("D53").Value = 1.427
$wsGainers.Range("E53").Value = 8.406000000000001
$wsGainers.Range("B54").Value = "MEDIASSIST"
$wsGainers.Range("C54").Value = 1.6755
$wsGainers.Range("D54").Value = 4.3257
$wsGainers.Range("E54").Value = 8.2271
$wsGainers.Range("B55").Value = "JKTYRE"
$wsGainers.Range("C55").Value = 1.6622
$wsGainers.Range("D55").Value = 4.6729
$wsGainers.Range("E55").Value = 20.5648
$wsGainers.Range("B56").Value = "HESTERBIO"
$wsGainers.Range("C56").Value = 1.639
$wsGainers.Range("D56").Value = -2.4748
$wsGainers.Range("E56").Value = -14.3267
$wsGainers.Range("B57").Value = "OIL"
$wsGainers.Range("C57").Value = 1.5814
$wsGainers.Range("D57").Value = 1.8236
$wsGainers.Range("E57").Value = 3.2262
$wsGainers.Range("B58").Value = "BIOCON"
$wsGainers.Range("C58").Value = 1.5788
$wsGainers.Range("D58").Value = 5.5911
$wsGainers.Range("E58").Value = 11.3196
$wsGainers.Range("B59").Value = "CREDITACC"
$wsGainers.Range("C59").Value = 1.553
$wsGainers.Range("D59").Value = 0.2109
$wsGainers.Range("E59").Value = 5.3429
$wsGainers.Range("B60").Value = "VIMTALABS"
$wsGainers.Range("C60").Value = 1.5524
$wsGainers.Range("D60").Value = 1.7262
$wsGainers.Range("E60").Value = -3.2413
$wsGainers.Range("B62").Value = "HGINFRA"
$wsGainers.Range("C62").Value = 1.5415
$wsGainers.Range("D62").Value = 2.9622
$wsGainers.Range("E62").Value = 1.3067
$wsGainers.Range("B63").Value = "MANGLMCEM"
$wsGainers.Range("C63").Value = 1.5329
$wsGainers.Range("D63").Value = -0.038
$wsGainers.Range("E63").Value = 5.5719
$wsGainers.Range("B64").Value = "SIMPLEXINF"
$wsGainers.Range("C64").Value = 1.5294
$wsGainers.Range("D64").Value = -1.9049
$wsGainers.Range("E64").Value = -2.0444
$wsGainers.Range("B65").Value = "ACUTAAS"
$wsGainers.Range("C65").Value = 1.5095
$wsGainers.Range("D65").Value = 1.8713
$wsGainers.Range("E65").Value = 30.126
$wsGainers.Range("B66").Value = "SOLEX"
$wsGainers.Range("C66").Value = 1.5014
$wsGainers.Range("D66").Value = 3.0204
$wsGainers.Range("E66").Value = "N/A"
$wsGainers.Range("B67").Value = "SUZLON"
$wsGainers.Range("C67").Value = 1.4951
$wsGainers.Range("D67").Value = 9.7362
$wsGainers.Range("E67").Value = 7.2843
$wsGainers.Range("B68").Value = "CRAMC"
$wsGainers.Range("C68").Value = 1.4412
$wsGainers.Range("D68").Value = 7.5053
$wsGainers.Range("E68").Value = "N/A"
$wsGainers.Range("B69").Value = "RAMCOIND"
$wsGainers.Range("C69").Value = 1.4336
$wsGainers.Range("D69").Value = 6.2105
$wsGainers.Range("E69").Value = 8.184699999999999
$wsGainers.Range("B70").Value = "RANEHOLDIN"
$wsGainers.Range("C70").Value = 1.4229
$wsGainers.Range("D70").Value = 2.705
$wsGainers.Range("E70").Value = -3.0841
$wsGainers.Range("B71").Value = "RTNPOWER"
$wsGainers.Range("C71").Value = 1.3974
$wsGainers.Range("D71").Value = 4.5005
$wsGainers.Range("E71").Value = 3.4759
$wsGainers.Range("B72").Value = "SUNTECK"
$wsGainers.Range("C72").Value = 1.3831
$wsGainers.Range("D72").Value = 3.6918
$wsGainers.Range("E72").Value = 8.587300000000001
$wsGainers.Range("B73").Value = "MUFIN"
$wsGainers.Range("C73").Value = 1.3796
$wsGainers.Range("D73").Value = 3.116
$wsGainers.Range("E73").Value = 13.5553
$wsGainers.Range("B74").Value = "RPTECH"
$wsGainers.Range("C74").Value = 1.3786
$wsGainers.Range("D74").Value = 3.1679
$wsGainers.Range("E74").Value = 5.5401
$wsGainers.Range("B75").Value = "EMBDL"
$wsGainers.Range("C75").Value = 1.3709
$wsGainers.Range("D75").Value = 1.414
$wsGainers.Range("E75").Value = 3.4262
$wsGainers.Range("C76").Value = 1.3617
$wsGainers.Range("D76").Value = 0.7165
$wsGainers.Range("E76").Value = 3.9749

# --- Top Losers sheet ---
$wsLosers = $wb.Worksheets.Item("Top Losers")
$wsLosers.Range("C2").Value = -15.5658
$wsLosers.Range("D2").Value = -14.2589
$wsLosers.Range("E2").Value = 3.3937
$wsLosers.Range("C3").Value = -5.9786
$wsLosers.Range("D3").Value = -2.4059
$wsLosers.Range("E3").Value = 10.3128
$wsLosers.Range("B4").Value = "STALLION"
$wsLosers.Range("C4").Value = -5.0008
$wsLosers.Range("D4").Value = -9.9626
$wsLosers.Range("E4").Value = 15.3662
$wsLosers.Range("B5").Value = "NSIL"
$wsLosers.Range("C5").Value = -4.9994
$wsLosers.Range("D5").Value = -6.6757
$wsLosers.Range("E5").Value = -0.4934
$wsLosers.Range("B6").Value = "KICL"
$wsLosers.Range("C6").Value = -4.999
$wsLosers.Range("D6").Value = -5.7975
$wsLosers.Range("E6").Value = 15.7082
$wsLosers.Range("B7").Value = "PILANIINVS"
$wsLosers.Range("C7").Value = -4.9969
$wsLosers.Range("D7").Value = -5.748
$wsLosers.Range("E7").Value = -0.9429999999999999
$wsLosers.Range("B8").Value = "DRREDDY"
$wsLosers.Range("C8").Value = -4.5567
$wsLosers.Range("D8").Value = -6.9882
$wsLosers.Range("E8").Value = -2.4352
$wsLosers.Range("B9").Value = "IDEA"
$wsLosers.Range("C9").Value = -4.3803
$wsLosers.Range("D9").Value = -6.9647
$wsLosers.Range("E9").Value = 10.0861
$wsLosers.Range("B10").Value = "DREDGECORP"
$wsLosers.Range("C10").Value = -3.5693
$wsLosers.Range("D10").Value = 17.4164
$wsLosers.Range("E10").Value = 18.184
$wsLosers.Range("B11").Value = "CGCL"
$wsLosers.Range("C11").Value = -3.3958
$wsLosers.Range("D11").Value = -1.4539
$wsLosers.Range("E11").Value = 9.0723
$wsLosers.Range("B12").Value = "TCI"
$wsLosers.Range("C12").Value = -3.0873
$wsLosers.Range("D12").Value = 0.2275
$wsLosers.Range("E12").Value = 0.7111
$wsLosers.Range("C13").Value = -3.0397
$wsLosers.Range("D13").Value = -1.7757
$wsLosers.Range("E13").Value = -2.5858
$wsLosers.Range("B14").Value = "MOLDTKPAC"
$wsLosers.Range("C14").Value = -2.94
$wsLosers.Range("D14").Value = -4.1495
$wsLosers.Range("E14").Value = -2.4877
$wsLosers.Range("B15").Value = "LXCHEM"
$wsLosers.Range("C15").Value = -2.8249
$wsLosers.Range("D15").Value = -3.2629
$wsLosers.Range("E15").Value = -4.4335
$wsLosers.Range("B16").Value = "INDUSTOWER"
$wsLosers.Range("C16").Value = -2.6899
$wsLosers.Range("D16").Value = 2.5584
$wsLosers.Range("E16").Value = 8.1365
$wsLosers.Range("B17").Value = "HCG"
$wsLosers.Range("C17").Value = -2.6729
$wsLosers.Range("D17").Value = -0.4677
$wsLosers.Range("E17").Value = 17.4491
$wsLosers.Range("B18").Value = "UBL"
$wsLosers.Range("C18").Value = -2.6409
$wsLosers.Range("D18").Value = -2.0447
$wsLosers.Range("E18").Value = -0.5053
$wsLosers.Range("B19").Value = "EPACKPEB"
$wsLosers.Range("C19").Value = -2.6121
$wsLosers.Range("D19").Value = -2.5468
$wsLosers.Range("E19").Value = "N/A"
$wsLosers.Range("B20").Value = "BHARATWIRE"
$wsLosers.Range("C20").Value = -2.3154
$wsLosers.Range("D20").Value = 19.9896
$wsLosers.Range("E20").Value = 21.0292
$wsLosers.Range("B21").Value = "INFOBEAN"
$wsLosers.Range("C21").Value = -2.2712
$wsLosers.Range("D21").Value = 20.2937
$wsLosers.Range("E21").Value = 35.0238
$wsLosers.Range("B22").Value = "MGL"
$wsLosers.Range("C22").Value = -2.1843
$wsLosers.Range("D22").Value = -2.8963
$wsLosers.Range("E22").Value = -2.1465
$wsLosers.Range("B23").Value = "KALAMANDIR"
$wsLosers.Range("C23").Value = -2.1663
$wsLosers.Range("D23").Value = -0.459
$wsLosers.Range("E23").Value = 23.2701
$wsLosers.Range("B24").Value = "SINDHUTRAD"
$wsLosers.Range("C24").Value = -2.1606
$wsLosers.Range("D24").Value = -1.0717
$wsLosers.Range("E24").Value = -15.0442
$wsLosers.Range("B25").Value = "BLUEDART"
$wsLosers.Range("C25").Value = -2.0085
$wsLosers.Range("D25").Value = 15.7129
$wsLosers.Range("E25").Value = 12.923
$wsLosers.Range("B26").Value = "STARHEALTH"
$wsLosers.Range("C26").Value = -1.9045
$wsLosers.Range("D26").Value = -3.432
$wsLosers.Range("E26").Value = 5.4952
$wsLosers.Range("B27").Value = "GODIGIT"
$wsLosers.Range("C27").Value = -1.901
$wsLosers.Range("D27").Value = 1.8169
$wsLosers.Range("E27").Value = 3.6411
$wsLosers.Range("C28").Value = -1.7949
$wsLosers.Range("D28").Value = 1.949
$wsLosers.Range("E28").Value = 5.766
$wsLosers.Range("B29").Value = "BHARTIHEXA"
$wsLosers.Range("C29").Value = -1.7866
$wsLosers.Range("D29").Value = 5.1745
$wsLosers.Range("E29").Value = 13.2727
$wsLosers.Range("B30").Value = "LICHSGFIN"
$wsLosers.Range("C30").Value = -1.7356
$wsLosers.Range("D30").Value = 0.6299
$wsLosers.Range("E30").Value = 3.2124
$wsLosers.Range("B31").Value = "INDOTHAI"
$wsLosers.Range("C31").Value = -1.7353
$wsLosers.Range("D31").Value = 2.7209
$wsLosers.Range("E31").Value = 41.2536
$wsLosers.Range("B32").Value = "MAHASTEEL"
$wsLosers.Range("C32").Value = -1.6393
$wsLosers.Range("D32").Value = 5.3186
$wsLosers.Range("E32").Value = 44.4565
$wsLosers.Range("B33").Value = "63MOONS"
$wsLosers.Range("C33").Value = -1.6353
$wsLosers.Range("D33").Value = 1.7812
$wsLosers.Range("E33").Value = -4.4833
$wsLosers.Range("B34").Value = "HMT"
$wsLosers.Range("C34").Value = -1.6337
$wsLosers.Range("D34").Value = -2.1269
$wsLosers.Range("E34").Value = -5.5407
$wsLosers.Range("B35").Value = "BHARTIARTL"
$wsLosers.Range("C35").Value = -1.5281
$wsLosers.Range("D35").Value = 1.9317
$wsLosers.Range("E35").Value = 10.1203
$wsLosers.Range("B36").Value = "HEG"
$wsLosers.Range("C36").Value = -1.5157
$wsLosers.Range("D36").Value = 10.8784
$wsLosers.Range("E36").Value = 13.2501
$wsLosers.Range("B37").Value = "SUNPHARMA"
$wsLosers.Range("C37").Value = -1.4336
$wsLosers.Range("D37").Value = -0.4532
$wsLosers.Range("E37").Value = 6.0842
$wsLosers.Range("B38").Value = "APOLLOPIPE"
$wsLosers.Range("C38").Value = -1.4037
$wsLosers.Range("D38").Value = -3.123
$wsLosers.Range("E38").Value = -8.2432
$wsLosers.Range("B39").Value = "NAM-INDIA"
$wsLosers.Range("C39").Value = -1.3347
$wsLosers.Range("D39").Value = -8.268800000000001
$wsLosers.Range("E39").Value = -2.522
$wsLosers.Range("B40").Value = "VEDL"
$wsLosers.Range("C40").Value = -1.2979
$wsLosers.Range("D40").Value = 2.8047
$wsLosers.Range("E40").Value = 9.3817
$wsLosers.Range("B41").Value = "YATRA"
$wsLosers.Range("C41").Value = -1.2609
$wsLosers.Range("D41").Value = -4.0705
$wsLosers.Range("E41").Value = 6.0172
$wsLosers.Range("B42").Value = "SGMART"
$wsLosers.Range("C42").Value = -1.2579
$wsLosers.Range("D42").Value = 6.8971
$wsLosers.Range("E42").Value = 1.2483
$wsLosers.Range("B43").Value = "ATLANTAELE"
$wsLosers.Range("C43").Value = -1.2549
$wsLosers.Range("D43").Value = -8.0814
$wsLosers.Range("E43").Value = 21.611
$wsLosers.Range("B44").Value = "AEGISLOG"
$wsLosers.Range("C44").Value = -1.2287
$wsLosers.Range("D44").Value = -0.8817
$wsLosers.Range("E44").Value = 1.2994
$wsLosers.Range("B45").Value = "SOLARWORLD"
$wsLosers.Range("C45").Value = -1.2264
$wsLosers.Range("D45").Value = 7.6836
$wsLosers.Range("E45").Value = 3.323
$wsLosers.Range("B46").Value = "KIOCL"
$wsLosers.Range("C46").Value = -1.1727
$wsLosers.Range("D46").Value = -4.3614
$wsLosers.Range("E46").Value = 1.2604
$wsLosers.Range("B47").Value = "ITC"
$wsLosers.Range("C47").Value = -1.1622
$wsLosers.Range("D47").Value = -0.024
$wsLosers.Range("E47").Value = 3.7729
$wsLosers.Range("B48").Value = "RMDRIP"
$wsLosers.Range("C48").Value = -1.1535
$wsLosers.Range("D48").Value = -0.0769
$wsLosers.Range("E48").Value = 2.8082
$wsLosers.Range("B49").Value = "VGUARD"
$wsLosers.Range("C49").Value = -1.1371
$wsLosers.Range("D49").Value = 1.612
$wsLosers.Range("E49").Value = 0.719
$wsLosers.Range("B50").Value = "JINDALPHOT"
$wsLosers.Range("C50").Value = -1.1154
$wsLosers.Range("D50").Value = -1.334
$wsLosers.Range("E50").Value = 21.6115
$wsLosers.Range("B51").Value = "SPAL"
$wsLosers.Range("C51").Value = -1.0797
$wsLosers.Range("D51").Value = 2.4269
$wsLosers.Range("E51").Value = 0.0207
$wsLosers.Range("B52").Value = "CUB"
$wsLosers.Range("C52").Value = -1.0715
$wsLosers.Range("D52").Value = 4.6691
$wsLosers.Range("E52").Value = 9.2921
$wsLosers.Range("B53").Value = "SURAJEST"
$wsLosers.Range("C53").Value = -1.0697
$wsLosers.Range("D53").Value = 8.0504
$wsLosers.Range("E53").Value = 6.0183
$wsLosers.Range("B54").Value = "FINOPB"
$wsLosers.Range("C54").Value = -1.0683
$wsLosers.Range("D54").Value = -7.1003
$wsLosers.Range("E54").Value = 10.2083
$wsLosers.Range("B55").Value = "NACLIND"
$wsLosers.Range("C55").Value = -1.0677
$wsLosers.Range("D55").Value = -3.006
$wsLosers.Range("E55").Value = 1.3752
$wsLosers.Range("B56").Value = "TARIL"
$wsLosers.Range("C56").Value = -1.0665
$wsLosers.Range("D56").Value = -5.499
$wsLosers.Range("E56").Value = -6.5673
$wsLosers.Range("B57").Value = "ABDL"
$wsLosers.Range("C57").Value = -1.0636
$wsLosers.Range("D57").Value = 3.8004
$wsLosers.Range("E57").Value = 26.4231
$wsLosers.Range("B58").Value = "GRWRHITECH"
$wsLosers.Range("C58").Value = -1.0595
$wsLosers.Range("D58").Value = -4.7011
$wsLosers.Range("E58").Value = 20.3878
$wsLosers.Range("B59").Value = "ZYDUSLIFE"
$wsLosers.Range("C59").Value = -1.042
$wsLosers.Range("D59").Value = -1.2389
$wsLosers.Range("E59").Value = 1.0744
$wsLosers.Range("B60").Value = "IDEAFORGE"
$wsLosers.Range("C60").Value = -1.0414
$wsLosers.Range("D60").Value = -0.2221
$wsLosers.Range("E60").Value = -2.003
$wsLosers.Range("B61").Value = "CPPLUS"
$wsLosers.Range("C61").Value = -1.0072
$wsLosers.Range("D61").Value = -1.6699
$wsLosers.Range("E61").Value = 3.5473
$wsLosers.Range("B62").Value = "APLAPOLLO"
$wsLosers.Range("C62").Value = -1.0029
$wsLosers.Range("D62").Value = 1.8411
$wsLosers.Range("E62").Value = 5.9664
$wsLosers.Range("B63").Value = "CMSINFO"
$wsLosers.Range("C63").Value = -1.0025
$wsLosers.Range("D63").Value = 1.6578
$wsLosers.Range("E63").Value = 1.862
$wsLosers.Range("B64").Value = "SAIL"
$wsLosers.Range("C64").Value = -0.9961
$wsLosers.Range("D64").Value = 7.4849
$wsLosers.Range("E64").Value = 3.4726
$wsLosers.Range("B65").Value = "INDIAMART"
$wsLosers.Range("C65").Value = -0.9905
$wsLosers.Range("D65").Value = 3.7245
$wsLosers.Range("E65").Value = 4.3584
$wsLosers.Range("B66").Value = "RAJRATAN"
$wsLosers.Range("C66").Value = -0.9875
$wsLosers.Range("D66").Value = -0.1764
$wsLosers.Range("E66").Value = 25.5508
$wsLosers.Range("B67").Value = "EUREKAFORB"
$wsLosers.Range("C67").Value = -0.9809
$wsLosers.Range("D67").Value = -0.6297
$wsLosers.Range("E67").Value = -1.7019
$wsLosers.Range("B68").Value = "ORIENTTECH"
$wsLosers.Range("C68").Value = -0.9801
$wsLosers.Range("D68").Value = -0.4605
$wsLosers.Range("E68").Value = 31.3781
$wsLosers.Range("B69").Value = "SAMBHV"
$wsLosers.Range("C69").Value = -0.9774
$wsLosers.Range("D69").Value = 1.621
$wsLosers.Range("E69").Value = 4.1391
$wsLosers.Range("B70").Value = "HFCL"
$wsLosers.Range("C70").Value = -0.9721
$wsLosers.Range("D70").Value = -1.7237
$wsLosers.Range("E70").Value = 4.8299
$wsLosers.Range("B71").Value = "DHANBANK"
$wsLosers.Range("C71").Value = -0.9691
$wsLosers.Range("D71").Value = 0.4537
$wsLosers.Range("E71").Value = 6.5784
$wsLosers.Range("B72").Value = "HDFCLIFE"
$wsLosers.Range("C72").Value = -0.9195
$wsLosers.Range("D72").Value = 2.6328
$wsLosers.Range("E72").Value = -0.2842
$wsLosers.Range("B73").Value = "PRECWIRE"
$wsLosers.Range("C73").Value = -0.9106
$wsLosers.Range("D73").Value = 11.3335
$wsLosers.Range("E73").Value = 22.1635
$wsLosers.Range("B74").Value = "AUROPHARMA"
$wsLosers.Range("C74").Value = -0.9085
$wsLosers.Range("D74").Value = 1.53
$wsLosers.Range("E74").Value = 1.633
$wsLosers.Range("B75").Value = "WEWORK"
$wsLosers.Range("C75").Value = -0.9036999999999999
$wsLosers.Range("D75").Value = 1.0555
$wsLosers.Range("E75").Value = "N/A"
$wsLosers.Range("B76").Value = "CANTABIL"
$wsLosers.Range("C76").Value = -0.9009
$wsLosers.Range("D76").Value = 3.9614
$wsLosers.Range("E76").Value = 3.2142

# --- 1 Month Performance sheet ---
$wsMonth = $wb.Worksheets.Item("1 Month Performance")
$wsMonth.Range("C2").Value = 109.2041
$wsMonth.Range("C4").Value = 79.863
$wsMonth.Range("C5").Value = 69.3258
$wsMonth.Range("C6").Value = 63.7151
$wsMonth.Range("B7").Value = "BGRENERGY"
$wsMonth.Range("C7").Value = 56.2986
$wsMonth.Range("B8").Value = "MAHASTEEL"
$wsMonth.Range("C8").Value = 53.4527
$wsMonth.Range("C9").Value = 49.729
$wsMonth.Range("C10").Value = 41.8348
$wsMonth.Range("B11").Value = "TVSELECT"
$wsMonth.Range("C11").Value = 40.7901
$wsMonth.Range("B12").Value = "NETWEB"
$wsMonth.Range("C12").Value = 40.5517
$wsMonth.Range("B13").Value = "STALLION"
$wsMonth.Range("C13").Value = 40.5383
$wsMonth.Range("B14").Value = "BHARATSE"
$wsMonth.Range("C14").Value = 39.4527
$wsMonth.Range("B15").Value = "RAMAPHO"
$wsMonth.Range("C15").Value = 38.5444
$wsMonth.Range("B16").Value = "TVSSRICHAK"
$wsMonth.Range("C16").Value = 38.1604
$wsMonth.Range("B17").Value = "SANDUMA"
$wsMonth.Range("C17").Value = 37.9959
$wsMonth.Range("B18").Value = "V2RETAIL"
$wsMonth.Range("C18").Value = 37.5703
$wsMonth.Range("C19").Value = 37.2425
$wsMonth.Range("B20").Value = "SHAREINDIA"
$wsMonth.Range("C20").Value = 37.2266
$wsMonth.Range("C21").Value = 36.1246
$wsMonth.Range("B22").Value = "MAANALU"
$wsMonth.Range("C22").Value = 35.531
$wsMonth.Range("B23").Value = "SOUTHBANK"
$wsMonth.Range("C23").Value = 34.7645
$wsMonth.Range("C24").Value = 30.4928
$wsMonth.Range("C25").Value = 29.1026
$wsMonth.Range("B26").Value = "RAMCOSYS"
$wsMonth.Range("C26").Value = 28.9085
$wsMonth.Range("C27").Value = 28.669
$wsMonth.Range("B29").Value = "ATHERENERG"
$wsMonth.Range("C29").Value = 28.4199
$wsMonth.Range("C31").Value = 27.9125
$wsMonth.Range("C32").Value = 26.6252
$wsMonth.Range("C34").Value = 26.3484
$wsMonth.Range("C35").Value = 25.7886
$wsMonth.Range("B36").Value = "SEJALLTD"
$wsMonth.Range("C36").Value = 25.7303
$wsMonth.Range("B37").Value = "AVALON"
$wsMonth.Range("C37").Value = 25.1919
$wsMonth.Range("B38").Value = "CPEDU"
$wsMonth.Range("C38").Value = 25.0334
$wsMonth.Range("B39").Value = "DCBBANK"
$wsMonth.Range("C39").Value = 24.8134
$wsMonth.Range("B40").Value = "MRPL"
$wsMonth.Range("C40").Value = 24.5014
$wsMonth.Range("B41").Value = "AUBANK"
$wsMonth.Range("C41").Value = 23.8189
$wsMonth.Range("B42").Value = "TDPOWERSYS"
$wsMonth.Range("C42").Value = 23.7154
$wsMonth.Range("C43").Value = 23.6814
$wsMonth.Range("B44").Value = "TATVA"
$wsMonth.Range("C44").Value = 23.5257
$wsMonth.Range("B45").Value = "CARTRADE"
$wsMonth.Range("C45").Value = 23.5202
$wsMonth.Range("B46").Value = "LORDSCHLO"
$wsMonth.Range("C46").Value = 22.8247
$wsMonth.Range("B47").Value = "SUBROS"
$wsMonth.Range("C47").Value = 22.4991
$wsMonth.Range("B48").Value = "SAGILITY"
$wsMonth.Range("C48").Value = 22.4234
$wsMonth.Range("B49").Value = "SURYODAY"
$wsMonth.Range("C49").Value = 22.1995
$wsMonth.Range("C50").Value = 22.1699
$wsMonth.Range("B51").Value = "GUJTHEM"
$wsMonth.Range("C51").Value = 21.9016
$wsMonth.Range("B52").Value = "PRIVISCL"
$wsMonth.Range("C52").Value = 21.867
$wsMonth.Range("B53").Value = "SKYGOLD"
$wsMonth.Range("C53").Value = 21.0699
$wsMonth.Range("B54").Value = "RBLBANK"
$wsMonth.Range("C54").Value = 21.0438
$wsMonth.Range("B55").Value = "PRECWIRE"
$wsMonth.Range("C55").Value = 20.9768
$wsMonth.Range("B56").Value = "MOLDTECH"
$wsMonth.Range("C56").Value = 20.8361
$wsMonth.Range("B57").Value = "HINDCOPPER"
$wsMonth.Range("C57").Value = 20.7462
$wsMonth.Range("B58").Value = "SKMEGGPROD"
$wsMonth.Range("C58").Value = 20.7435
$wsMonth.Range("B59").Value = "IIFL"
$wsMonth.Range("C59").Value = 20.5731
$wsMonth.Range("B60").Value = "BHARATWIRE"
$wsMonth.Range("C60").Value = 20.5142
$wsMonth.Range("B61").Value = "BANKINDIA"
$wsMonth.Range("C61").Value = 20.4714
$wsMonth.Range("C62").Value = 19.8852
$wsMonth.Range("B63").Value = "SHRIRAMFIN"
$wsMonth.Range("C63").Value = 19.7784
$wsMonth.Range("B64").Value = "FEDERALBNK"
$wsMonth.Range("C64").Value = 19.7291
$wsMonth.Range("B65").Value = "ORBTEXP"
$wsMonth.Range("C65").Value = 19.6363
$wsMonth.Range("B66").Value = "ETHOSLTD"
$wsMonth.Range("C66").Value = 19.3771
$wsMonth.Range("B67").Value = "MCX"
$wsMonth.Range("C67").Value = 19.3513
$wsMonth.Range("B69").Value = "TERASOFT"
$wsMonth.Range("C69").Value = 19.2961
$wsMonth.Range("B70").Value = "GRMOVER"
$wsMonth.Range("C70").Value = 18.6571
$wsMonth.Range("B71").Value = "THOMASCOTT"
$wsMonth.Range("C71").Value = 18.6189
$wsMonth.Range("B72").Value = "LUMAXIND"
$wsMonth.Range("C72").Value = 18.4254
$wsMonth.Range("C73").Value = 18.4106
$wsMonth.Range("C74").Value = 18.3742
$wsMonth.Range("B75").Value = "TINNARUBR"
$wsMonth.Range("C75").Value = 18.1169
$wsMonth.Range("B76").Value = "BAJAJCON"
$wsMonth.Range("C76").Value = 18.0898

# --- distance from Dma50 sheet ---
$wsDma50 = $wb.Worksheets.Item("distance from Dma50")
$wsDma50.Range("C2").Value = 10.1337
$wsDma50.Range("C3").Value = 7.8245
$wsDma50.Range("C4").Value = 5.9022
$wsDma50.Range("C5").Value = 5.7561
$wsDma50.Range("C6").Value = 5.3193
$wsDma50.Range("C7").Value = 5.2003
$wsDma50.Range("C8").Value = 4.685
$wsDma50.Range("C9").Value = 4.6531
$wsDma50.Range("C10").Value = 3.8951
$wsDma50.Range("C11").Value = 3.732
$wsDma50.Range("C12").Value = 3.5654
$wsDma50.Range("C13").Value = 3.5142
$wsDma50.Range("C14").Value = 3.3129
$wsDma50.Range("C15").Value = 3.2793
$wsDma50.Range("C16").Value = 3.21
$wsDma50.Range("C17").Value = 3.0607
$wsDma50.Range("C18").Value = 2.9181
$wsDma50.Range("C19").Value = 2.6433
$wsDma50.Range("C20").Value = 2.6011
$wsDma50.Range("C21").Value = 2.4439
$wsDma50.Range("C22").Value = 1.986
$wsDma50.Range("C23").Value = 1.5207
$wsDma50.Range("C24").Value = 1.343
$wsDma50.Range("B25").Value = "NIFTYHEALTHCARE"
$wsDma50.Range("C25").Value = 1.2104
$wsDma50.Range("B26").Value = "NIFTYGROWSECT15"
$wsDma50.Range("C26").Value = 1.1852
$wsDma50.Range("C27").Value = 0.7163
$wsDma50.Range("C28").Value = 0.586
$wsDma50.Range("C29").Value = 0.093
$wsDma50.Range("C30").Value = -1.7542
